$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TipoIdentificacion")

# --- Add the three new rows of data (6: CC, 7: Pasaporte, 8: TI) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "CC"
$ws.Range("C7").Formula = "=B7"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Pasaporte"
$ws.Range("C8").Formula = "=B8"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "TI"
$ws.Range("C9").Formula = "=B9"

# --- Match formatting of the new rows to the existing table rows above them ---
# Columns A:B pick up the bordered style used by the rest of the table, then
# re-assert "no fill" explicitly (mirrors the extra applyFill flag Excel writes
# when the format is reconfirmed through the UI).
$ws.Range("C2").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)
$ws.Range("A7:B7").Interior.ColorIndex = -4142

$ws.Range("A7:B7").Copy()
$ws.Range("A8:B9").PasteSpecial(-4122)

# Column C keeps the highlighted "result" style used by the rest of the table.
$ws.Range("C2").Copy()
$ws.Range("C7:C9").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- Move the active selection to match where the user ended up editing ---
$ws.Range("I17").Select()
